# Duplicate the only worksheet ("Sheet1"). Excel's default "Move or Copy...
# (Create a copy)" behaviour inserts the new copy immediately before the
# sheet it was copied from and makes the copy the active sheet, naming it
# "Sheet1 (2)".
$wb = $excel.ActiveWorkbook
$original = $wb.Worksheets.Item(1)
$original.Copy($original) | Out-Null

# The copy is now the first/active sheet.
$copy = $wb.Worksheets.Item(1)

# Record the autograder breakdown for proj2 on the new copy only — the
# original sheet is left untouched.
$copy.Range("L6").Value = "31/31 code"
$copy.Range("M6").Value = "10/10 written"

# Leave the selection on the new sheet at M7, as it was after entering the
# note above.
$copy.Range("M7").Select() | Out-Null
